$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(
45810.01041666666,  45810.02083333334,  45810.03125,  45810.04166666666,  45810.05208333334,  45810.0625,  45810.07291666666,  45810.08333333334,  45810.09375,  45810.10416666666,  45810.11458333334,  45810.125,  45810.13541666666,  45810.14583333334,  45810.15625,  45810.16666666666,  45810.17708333334,  45810.1875,  45810.19791666666,  45810.20833333334,  45810.21875,  45810.22916666666,  45810.23958333334,  45810.25,  45810.26041666666,  45810.27083333334,  45810.28125,  45810.29166666666,  45810.30208333334,  45810.3125,  45810.32291666666,  45810.33333333334,  45810.34375,  45810.35416666666,  45810.36458333334,  45810.375,  45810.38541666666,  45810.39583333334,  45810.40625,  45810.41666666666,  45810.42708333334,  45810.4375,  45810.44791666666,  45810.45833333334,  45810.46875,  45810.47916666666,  45810.48958333334,  45810.5,  45810.51041666666,  45810.52083333334,  45810.53125,  45810.54166666666,  45810.55208333334,  45810.5625,  45810.57291666666,  45810.58333333334,  45810.59375,  45810.60416666666,  45810.61458333334,  45810.625,  45810.63541666666,  45810.64583333334,  45810.65625,  45810.66666666666,  45810.67708333334,  45810.6875,  45810.69791666666,  45810.70833333334,  45810.71875,  45810.72916666666,  45810.73958333334,  45810.75,  45810.76041666666,  45810.77083333334,  45810.78125,  45810.79166666666,  45810.80208333334,  45810.8125,  45810.82291666666,  45810.83333333334,  45810.84375,  45810.85416666666,  45810.86458333334,  45810.875,  45810.88541666666,  45810.89583333334,  45810.90625,  45810.91666666666,  45810.92708333334,  45810.9375,  45810.94791666666,  45810.95833333334,  45810.96875,  45810.97916666666,  45810.98958333334,  45811
)

$bValues = @(
95,  100,  100,  101,  113,  115,  117,  120,  148,  148,  149,  150,  167,  163,  164,  165,  159,  160,  160,  155,  150,  147,  146,  146,  105,  106,  106,  107,  75,  76,  77,  74,  57,  59,  61,  57,  47,  47,  47,  47,  45,  45,  45,  46,  54,  55,  55,  55,  79,  80,  81,  81,  108,  108,  109,  109,  141,  142,  142,  142,  162,  162,  162,  162,  178,  178,  178,  179,  201,  202,  203,  204,  250,  251,  253,  256,  315,  318,  322,  325,  379,  381,  383,  385,  422,  423,  425,  427,  436,  437,  438,  439,  0,  0,  0,  0
)

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}
